$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Strip the leading manual numbering ("1 . ", "2. ", "3.", "4. ",
#    "5. ") from the paragraphs that will become an automatic
#    numbered list, and make small wording fixes from the diff.
# ---------------------------------------------------------------

$d.Content.Find.Execute(
    "1 . Se ingresa un texto de varias",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Se ingresa un texto de varias", 2) | Out-Null

$d.Content.Find.Execute(
    "2. Se ingresa un texto o frase",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Se ingresa un texto o frase", 2) | Out-Null

$d.Content.Find.Execute(
    "3.Teniendo en cuenta el ejercicio 2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Teniendo en cuenta el ejercicio 2", 2) | Out-Null

$d.Content.Find.Execute(
    "4. Teniendo en cuenta el ejercicio 3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Teniendo en cuenta el ejercicio 3", 2) | Out-Null

$d.Content.Find.Execute(
    "5. Leer un archivo de texto",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Leer un archivo de texto", 2) | Out-Null

# Fix the double space in "empleados  de" -> "empleados de"
$d.Content.Find.Execute(
    "Mostrar los empleados  de la empresa",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mostrar los empleados de la empresa", 2) | Out-Null

# ---------------------------------------------------------------
# 2) Append the punctuation the diff adds at the end of several
#    items (periods / colon).
# ---------------------------------------------------------------

function Append-Text($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $tail = $d.Range($r.End - 1, $r.End - 1)
    $tail.InsertAfter($text)
}

Append-Text 3 "."
Append-Text 4 "."
Append-Text 7 ":"
Append-Text 8 "."
Append-Text 9 "."
Append-Text 10 "."
Append-Text 12 "."

# ---------------------------------------------------------------
# 3) Turn paragraphs 3-12 into the "Prrafodelista" (List Paragraph)
#    numbered list, sharing one numId, alternating outline levels.
# ---------------------------------------------------------------

$listRange = $d.Range($d.Paragraphs(3).Range.Start, $d.Paragraphs(12).Range.End)
$listRange.Style = "Prrafodelista"
$listRange.ListFormat.ApplyNumberDefault()

# Demote the sub-items to the second level (ilvl = 1)
$d.Paragraphs(6).Range.ListFormat.ListIndent()
$d.Paragraphs(8).Range.ListFormat.ListIndent()
$d.Paragraphs(9).Range.ListFormat.ListIndent()
$d.Paragraphs(10).Range.ListFormat.ListIndent()
$d.Paragraphs(11).Range.ListFormat.ListIndent()
$d.Paragraphs(12).Range.ListFormat.ListIndent()
